{"js": "// Update the court placeholder in the document header table from\n// \"<<court>>\" to \"<<courtName>>\" so the template uses the court/site\n// name merge field instead of the old \"court\" field.\n//\n// The placeholder lives in the first cell of the first table, as the\n// paragraph \"In the County Court at <<court>>\". We locate the literal\n// run of text \"court\" (case-sensitive, so we don't match the word\n// inside \"County\"/\"Court\") within that cell and insert \"Name\"\n// immediately after it, turning \"<<court>>\" into \"<<courtName>>\"\n// while keeping the existing bold Arial formatting.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"Expected at least one table in the document.\");\n}\n\n// The header table is the first table in the document; its first cell\n// holds \"In the County Court at <<court>>\".\nconst headerTable = tables.items[0];\nconst headerCell = headerTable.getCell(0, 0);\nconst headerParagraph = headerCell.body.paragraphs.getFirst();\n\n// Search only within that paragraph for the exact \"court\" run so we\n// never touch the other \"<<court>>\" placeholder further down in the\n// document body.\nconst searchResults = headerParagraph.search(\"court\", { matchCase: true });\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error('Could not find \"court\" placeholder text to update.');\n}\n\n// There should be exactly one \"court\" occurrence in this paragraph;\n// use the first match.\nconst courtRange = searchResults.items[0];\ncourtRange.insertText(\"Name\", Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "# Update the court placeholder in the document header table from\n# \"<<court>>\" to \"<<courtName>>\" so the template uses the court/site\n# name merge field instead of the old \"court\" field.\n#\n# The placeholder lives in the first cell of the first table, as the\n# paragraph \"In the County Court at <<court>>\". We restrict our search\n# to that single table cell and look for the literal (case-sensitive)\n# run of text \"court\" so we never touch the other \"<<court>>\"\n# placeholder that appears later in the document body.\n\n$d = $word.ActiveDocument\n\n$headerTable = $d.Tables.Item(1)\n$headerCell = $headerTable.Cell(1, 1)\n$cellRange = $headerCell.Range\n\n$find = $cellRange.Find\n$find.ClearFormatting()\n$find.Text = \"court\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n$find.Wrap = 0  # wdFindStop - do not search past the cell range\n\n$found = $find.Execute()\nif (-not $found) {\n  throw \"Could not find the 'court' placeholder in the header table cell.\"\n}\n\n# $cellRange now spans just the matched \"court\" text; insert the\n# literal \"Name\" immediately after it so \"<<court>>\" becomes\n# \"<<courtName>>\", keeping the existing bold Arial formatting.\n$cellRange.InsertAfter(\"Name\")\n\n$finalText = $d.Tables.Item(1).Cell(1, 1).Range.Text\nWrite-Output \"Header cell text is now: $finalText\"\n"}
